# Update the "想去人数" (want-to-go count) figures in column F for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets to reflect freshly
# scraped numbers, as published by the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> (old, new) value updates shared between the two sheets, keyed by
# the row number within each sheet where the corresponding record lives.
$exhibitionUpdates = @{
    2  = 798
    5  = 1031
    9  = 372
    11 = 489
    12 = 524
    14 = 12290
    15 = 72
    16 = 5464
}

$allTypesUpdates = @{
    2  = 798
    7  = 1031
    11 = 372
    13 = 489
    14 = 524
    16 = 12290
    18 = 72
    19 = 5464
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
